$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("PDL Date") contains the value "04012025" for every data row
# (rows 2 through 394). Replace it with "99999999" for all those rows,
# keeping the header in B1 untouched. Force text format so the numeric
# looking string is not converted to a number by Excel.
$lastRow = 394
$rng = $ws.Range("B2:B$lastRow")
$rng.NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "99999999"
}
